$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(23, 1).Value = "Neurological/SOP -  Femoral site care.pdf"
$ws.Cells.Item(24, 1).Value = "Airway/Tracheostomy_Laryngectomy/Tracheostomy change in Critical Care.pdf"
$ws.Cells.Item(25, 1).Value = "Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf"
$ws.Cells.Item(29, 1).Value = "Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf"
$ws.Cells.Item(30, 1).Value = "Drugs/heparin for Haemofiltration.pdf"
$ws.Cells.Item(31, 1).Value = "Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf"
$ws.Cells.Item(33, 1).Value = "Covid-19/WGH/CoVid intubation checklist WGH.pdf"
$ws.Cells.Item(34, 1).Value = "Airway/Emergency intubation checklist_em_pub.pdf"
$ws.Cells.Item(41, 1).Value = "Delirium/Managing a Potentially Violent Patient.pdf"
$ws.Cells.Item(42, 1).Value = "Delirium/Risk assessment posi mit.pdf"
$ws.Cells.Item(43, 1).Value = "Infection_and_sepsis/SOP Ultrasound Cleaning.pdf"
$ws.Cells.Item(44, 1).Value = "GI_Liver_and_Transplant/Treatment of constipation.pdf"
$ws.Cells.Item(45, 1).Value = "GI_Liver_and_Transplant/Abdominal pressure measurement.pdf"
$ws.Cells.Item(46, 1).Value = "Airway/Anticipated difficult airway tool.pdf"
$ws.Cells.Item(47, 1).Value = "Drugs/ketamine_in_asthma.pdf"
$ws.Cells.Item(48, 1).Value = "Breathing(Respiratory)/HFNO.pdf"
$ws.Cells.Item(50, 1).Value = "Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf"
$ws.Cells.Item(51, 1).Value = "Delirium/Drugs Causing Delirium and Agitiation.pdf"
$ws.Cells.Item(52, 1).Value = "Neurological/Sub arachnoid haemorrhage management.pdf"
$ws.Cells.Item(59, 1).Value = "Cardiovascular/Steroids for Septic Shock.pdf"
$ws.Cells.Item(60, 1).Value = "Breathing(Respiratory)/Equipment/APRV.pdf"
$ws.Cells.Item(62, 1).Value = "Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf"
$ws.Cells.Item(63, 1).Value = "Breathing(Respiratory)/Equipment/T piece Y piece.pdf"
$ws.Cells.Item(71, 1).Value = "Drugs/Antibiotic doses in CVVHD.pdf"
$ws.Cells.Item(72, 1).Value = "Diabetes_and_Glucose/Intravenous Insulin Therapy (not for DKA or HHS).pdf"
$ws.Cells.Item(73, 1).Value = "GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf"
$ws.Cells.Item(74, 1).Value = "GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf"
$ws.Cells.Item(76, 1).Value = "Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf"
$ws.Cells.Item(78, 1).Value = "Drugs/vasopressin organ donation.pdf"
$ws.Cells.Item(81, 1).Value = "Breathing(Respiratory)/Equipment/Bipap V60.pdf"
$ws.Cells.Item(82, 1).Value = "Breathing(Respiratory)/CPAP.pdf"
$ws.Cells.Item(83, 1).Value = "Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf"
$ws.Cells.Item(84, 1).Value = "Infection_and_sepsis/Infection indications for IVIG.pdf"
$ws.Cells.Item(85, 1).Value = "Drugs/piperacillin_tazobactam extended_infusion.pdf"
$ws.Cells.Item(88, 1).Value = "Neurological/Treatment of status epilepticus.pdf"
$ws.Cells.Item(89, 1).Value = "Routine_Care/Video Communication.pdf"
$ws.Cells.Item(93, 1).Value = "Cardiovascular/Management of hypertension within Critical Care.pdf"
$ws.Cells.Item(94, 1).Value = "Drugs/aminophylline.pdf"
$ws.Cells.Item(95, 1).Value = "Haematology_CAR-T/CRS.pdf"
$ws.Cells.Item(97, 1).Value = "Haematology_CAR-T/ICANS.pdf"
$ws.Cells.Item(98, 1).Value = "Drugs/phenytoin.pdf"
$ws.Cells.Item(99, 1).Value = "Drugs/rocuronium.pdf"
$ws.Cells.Item(100, 1).Value = "Drugs/milrinone.pdf"
$ws.Cells.Item(101, 1).Value = "Policies_and_admin/General Critical Care SOP_pub.pdf"
$ws.Cells.Item(102, 1).Value = "GI_Liver_and_Transplant/ICU - Upper GI bleeding (Endoscopy guideline).pdf"
$ws.Cells.Item(103, 1).Value = "Drugs/clonidine.pdf"
$ws.Cells.Item(104, 1).Value = "Neurological/Critical Care MRI Procedure_pub.pdf"
$ws.Cells.Item(105, 1).Value = "Ethics_and_Law/DNACPR policy for Scotland.pdf"
$ws.Cells.Item(106, 1).Value = "End_of_life_care/Palliative extubation & withdrawal of invasive ventilatory support nursing checklist.pdf"
$ws.Cells.Item(107, 1).Value = "End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD.pdf"
$ws.Cells.Item(108, 1).Value = "Neurological/Management of traumatic brain injury.pdf"
$ws.Cells.Item(109, 1).Value = "Organ_donation/Organ Retrieval SOP.pdf"
$ws.Cells.Item(111, 1).Value = "Infection_and_sepsis/Trip Out of Unit infection guidance.pdf"
$ws.Cells.Item(112, 1).Value = "Ethics_and_Law/Care at the End of Life (FICM).pdf"
$ws.Cells.Item(117, 1).Value = "GI_Liver_and_Transplant/Confirmation of Nasogastric Tube Position.pdf"
$ws.Cells.Item(120, 1).Value = "Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf"
$ws.Cells.Item(126, 1).Value = "Drugs/dobutamine.pdf"
$ws.Cells.Item(127, 1).Value = "Drugs/adrenaline.pdf"
$ws.Cells.Item(130, 1).Value = "Drugs/Alteplase for massive PE.pdf"
$ws.Cells.Item(131, 1).Value = "Drugs/alfentanil.pdf"
$ws.Cells.Item(133, 1).Value = "Drugs/all IV drug infusion information.pdf"
$ws.Cells.Item(134, 1).Value = "Drugs/neostigmine.pdf"
$ws.Cells.Item(135, 1).Value = "Drugs/vancomycin.pdf"
$ws.Cells.Item(136, 1).Value = "Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf"
$ws.Cells.Item(137, 1).Value = "Drugs/labetalol.pdf"
$ws.Cells.Item(139, 1).Value = "Drugs/midazolam.pdf"
$ws.Cells.Item(140, 1).Value = "Cardiovascular/Management of Acute Type B Aortic Dissection Guideline.pdf"
$ws.Cells.Item(142, 1).Value = "Procedures/CVC Guidance/CVC NHL  April 2023.pdf"
$ws.Cells.Item(143, 1).Value = "Drugs/salbutamol.pdf"
$ws.Cells.Item(145, 1).Value = "Routine_Care/ICU Eye Care Guideline.pdf"
$ws.Cells.Item(146, 1).Value = "Drugs/amiodarone.pdf"
$ws.Cells.Item(147, 1).Value = "Drugs/nicardipine.pdf"
$ws.Cells.Item(148, 1).Value = "Drugs/phenobarbitone.pdf"
$ws.Cells.Item(149, 1).Value = "Procedures/Arterial Line insertion for ACCPs.pdf"
$ws.Cells.Item(150, 1).Value = "Breathing(Respiratory)/Manual Ventilation and MHI.pdf"
$ws.Cells.Item(151, 1).Value = "Drugs/noradrenaline (peripheral).pdf"
$ws.Cells.Item(152, 1).Value = "Neurological/Ventriculitis Guideline.pdf"
$ws.Cells.Item(154, 1).Value = "Drugs/Epoprostenol.pdf"
$ws.Cells.Item(156, 1).Value = "Cardiovascular/Cardiac Output Monitoring_pub .pdf"
$ws.Cells.Item(157, 1).Value = "Cardiovascular/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf"
$ws.Cells.Item(158, 1).Value = "Drugs/calcium.pdf"
$ws.Cells.Item(159, 1).Value = "Drugs/dalteparin_thromboprophylaxis.pdf"
$ws.Cells.Item(160, 1).Value = "Drugs/Vancomycin Continuous Infusion Fluid Restricted.pdf"
$ws.Cells.Item(162, 1).Value = "Airway/Tracheostomy_Laryngectomy/Decannulation Guidline.pdf"
$ws.Cells.Item(167, 1).Value = "Drugs/ketamine_for_status epilepticus.pdf"
$ws.Cells.Item(168, 1).Value = "Drugs/Phosphate.pdf"
$ws.Cells.Item(169, 1).Value = "Drugs/Thiopentone.pdf"
$ws.Cells.Item(170, 1).Value = "Breathing(Respiratory)/Proning Guideline.pdf"
$ws.Cells.Item(172, 1).Value = "Drugs/Octreotide.pdf"
$ws.Cells.Item(173, 1).Value = "Procedures/ACCP CVC placement following completion of initial competencies.pdf"
$ws.Cells.Item(174, 1).Value = "Procedures/ACCPs acquiring initial CVC competencies.pdf"
$ws.Cells.Item(175, 1).Value = "Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf"
$ws.Cells.Item(179, 1).Value = "End_of_life_care/Guideline following Sudden Cardiac Death where death occurs in ICU.pdf"
$ws.Cells.Item(180, 1).Value = "Breathing(Respiratory)/Equipment/NIV through Drager Vent Set up in Critical Care.pdf"
